$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 33.715737664422008
$ws.Range("B3").Value = 2.8643075125809778
$ws.Range("B4").Value = 17.672946448842847
$ws.Range("B5").Value = 46.868462385135516
$ws.Range("B6").Value = 10.836205237344824
$ws.Range("B7").Value = 17.843623424291518
$ws.Range("B8").Value = 30.488507146666773
$ws.Range("B9").Value = 15.695620009799692
